$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'244.85"
$ws.Range('E2').Value = "'-0.52%"
$ws.Range('G2').Value = "'23"
$ws.Range('D3').Value = "'27.34"
$ws.Range('E3').Value = "'4.89%"
$ws.Range('G3').Value = "'23"
$ws.Range('D4').Value = "'5.112"
$ws.Range('E4').Value = "'0.44%"
$ws.Range('G4').Value = "'23"
$ws.Range('D5').Value = "'0.05678"
$ws.Range('E5').Value = "'1.71%"
$ws.Range('G5').Value = "'23"
$ws.Range('D6').Value = "'6.492"
$ws.Range('E6').Value = "'0.07%"
$ws.Range('G6').Value = "'23"
$ws.Range('E7').Value = "'0.76%"
$ws.Range('G7').Value = "'23"
$ws.Range('D8').Value = "'0.8497"
$ws.Range('E8').Value = "'0.77%"
$ws.Range('G8').Value = "'23"
$ws.Range('B9').Value = 'MandalaExchangeToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D9').Value = "'0.06956"
$ws.Range('E9').Value = "'0.44%"
$ws.Range('G9').Value = "'23"
$ws.Range('B10').Value = 'BitrueCoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D10').Value = "'0.02879"
$ws.Range('E10').Value = "'2.26%"
$ws.Range('G10').Value = "'23"
$ws.Range('B11').Value = 'BitMartToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D11').Value = "'0.09392"
$ws.Range('E11').Value = "'0.23%"
$ws.Range('G11').Value = "'23"
$ws.Range('B12').Value = 'BitForexToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D12').Value = "'0.001509"
$ws.Range('E12').Value = "'-0.91%"
$ws.Range('G12').Value = "'23"
$ws.Range('B13').Value = 'CoinExToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D13').Value = "'0.04028"
$ws.Range('E13').Value = "'-13.73%"
$ws.Range('G13').Value = "'23"
$ws.Range('B14').Value = 'One'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').Value = "'0.0005981"
$ws.Range('E14').Value = "'0.30%"
$ws.Range('G14').Value = "'23"
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = "'0.006215"
$ws.Range('E15').Value = "'0.20%"
$ws.Range('G15').Value = "'23"
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = "'3.512"
$ws.Range('E16').Value = "'-2.61%"
$ws.Range('G16').Value = "'23"
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = "'3.010"
$ws.Range('E17').Value = "'-0.39%"
$ws.Range('G17').Value = "'23"
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = "'2.227"
$ws.Range('E18').Value = "'8.36%"
$ws.Range('G18').Value = "'23"
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = "'0.3158"
$ws.Range('E19').Value = "'1.48%"
$ws.Range('G19').Value = "'23"
$ws.Range('B20').Value = 'WazirX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D20').Value = "'0.1334"
$ws.Range('E20').Value = "'0.00%"
$ws.Range('G20').Value = "'23"
$ws.Range('G21').Value = "'23"
$ws.Range('E22').Value = "'-1.57%"
$ws.Range('G22').Value = "'23"
$ws.Range('D23').Value = "'3.557"
$ws.Range('E23').Value = "'-5.35%"
$ws.Range('G23').Value = "'23"
$ws.Range('E24').Value = "'-0.05%"
$ws.Range('G24').Value = "'23"
$ws.Range('D25').Value = "'0.001217"
$ws.Range('E25').Value = "'-2.67%"
$ws.Range('G25').Value = "'23"
$ws.Range('D26').Value = "'0.004472"
$ws.Range('E26').Value = "'-1.79%"
$ws.Range('G26').Value = "'23"
$ws.Range('E27').Value = "'22.82%"
$ws.Range('G27').Value = "'23"
$ws.Range('E28').Value = "'-27.46%"
$ws.Range('G28').Value = "'23"
$ws.Range('G29').Value = "'23"
$ws.Range('G30').Value = "'23"
$ws.Range('G31').Value = "'23"
$ws.Range('G32').Value = "'23"
$ws.Range('G33').Value = "'23"
$ws.Range('G34').Value = "'23"
$ws.Range('G35').Value = "'23"
$ws.Range('G36').Value = "'23"
$ws.Range('G37').Value = "'23"
$ws.Range('G38').Value = "'23"
$ws.Range('G39').Value = "'23"
$ws.Range('D40').Value = "'0.03716"
$ws.Range('E40').Value = "'1.80%"
$ws.Range('G40').Value = "'23"
$ws.Range('D41').Value = "'0.005950"
$ws.Range('E41').Value = "'-3.33%"
$ws.Range('G41').Value = "'23"
$ws.Range('E42').Value = "'0.44%"
$ws.Range('G42').Value = "'23"
$ws.Range('D43').Value = "'0.002414"
$ws.Range('E43').Value = "'-7.30%"
$ws.Range('G43').Value = "'23"
$ws.Range('D44').Value = "'0.009713"
$ws.Range('E44').Value = "'21.64%"
$ws.Range('G44').Value = "'23"
$ws.Range('D45').Value = "'0.00005102"
$ws.Range('E45').Value = "'-5.62%"
$ws.Range('G45').Value = "'23"
$ws.Range('E46').Value = "'-0.07%"
$ws.Range('G46').Value = "'23"
$ws.Range('E47').Value = "'-30.37%"
$ws.Range('G47').Value = "'23"
$ws.Range('D48').Value = "'0.002528"
$ws.Range('E48').Value = "'5.15%"
$ws.Range('G48').Value = "'23"
$ws.Range('E49').Value = "'-0.07%"
$ws.Range('G49').Value = "'23"
$ws.Range('E50').Value = "'-0.07%"
$ws.Range('G50').Value = "'23"
$ws.Range('G51').Value = "'23"
